$d = $word.ActiveDocument

$replacements = @(
    @("431÷4=", "453÷7="),
    @("572÷4=", "289÷8="),
    @("437÷8=", "927÷3="),
    @("408÷7=", "380÷5="),
    @("382÷4=", "653÷2="),
    @("957÷9=", "556÷4="),
    @("508÷2=", "963÷9="),
    @("253÷4=", "135÷4="),
    @("577÷3=", "172÷3="),
    @("234÷6=", "420÷8="),
    @("462÷9=", "260÷3="),
    @("878÷3=", "341÷9="),
    @("154÷4=", "878÷4="),
    @("926÷8=", "130÷3="),
    @("883÷6=", "940÷9="),
    @("514÷2=", "446÷2="),
    @("523÷7=", "963÷9="),
    @("453÷8=", "685÷7="),
    @("454÷7=", "769÷3="),
    @("297÷2=", "380÷2="),
    @("815÷6=", "899÷8="),
    @("755÷3=", "125÷3="),
    @("956÷2=", "612÷2="),
    @("950÷2=", "494÷8="),
    @("770÷9=", "180÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
